# "break out stock.yaml completed"
#
# 1) On the "day" sheet, six bsecode cells (D658:D663) were being written
#    as text ("532827" etc.) instead of numbers. Fix their type to numeric,
#    keeping the same values.
#
# 2) On the "week" sheet, a fresh scrape ("04/10/2024 11:32:49") appended
#    twelve new rows (374-385) after the existing data (which ended at 373).

$wb  = $excel.ActiveWorkbook
$day = $wb.Worksheets.Item("day")
$week = $wb.Worksheets.Item("week")

# --- 1) day!D658:D663 : text "532827" etc -> numeric 532827 etc ---------
$day.Range("D658").Value = 532827
$day.Range("D659").Value = 500488
$day.Range("D660").Value = 533179
$day.Range("D661").Value = 539524
$day.Range("D662").Value = 500331
$day.Range("D663").Value = 532830

# --- 2) week!A374:I385 : append the newly scraped rows -------------------
# bsecode (column D) keeps arriving as text from the scraper, so force the
# cells to Text format before writing the numeric-looking codes, otherwise
# they'd be auto-coerced to numbers.
$week.Range("D374:D385").NumberFormat = "@"

$week.Range("A374").Value = 1
$week.Range("B374").Value = "ABBOTINDIA"
$week.Range("C374").Value = "Abbott India Limited"
$week.Range("D374").Value = "500488"
$week.Range("E374").Value = -0.41
$week.Range("F374").Value = 28150
$week.Range("G374").Value = 10751
$week.Range("H374").Value = "week"
$week.Range("I374").Value = "04/10/2024 11:32:49"

$week.Range("A375").Value = 2
$week.Range("B375").Value = "LTTS"
$week.Range("C375").Value = "L&t Technology Services Limited"
$week.Range("D375").Value = "540115"
$week.Range("E375").Value = -0.5600000000000001
$week.Range("F375").Value = 5071.6
$week.Range("G375").Value = 319507
$week.Range("H375").Value = "week"
$week.Range("I375").Value = "04/10/2024 11:32:49"

$week.Range("A376").Value = 3
$week.Range("B376").Value = "TCS"
$week.Range("C376").Value = "Tata Consultancy Services Limited"
$week.Range("D376").Value = "532540"
$week.Range("E376").Value = 0.46
$week.Range("F376").Value = 4252.25
$week.Range("G376").Value = 2965470
$week.Range("H376").Value = "week"
$week.Range("I376").Value = "04/10/2024 11:32:49"

$week.Range("A377").Value = 4
$week.Range("B377").Value = "CUMMINSIND"
$week.Range("C377").Value = "Cummins India Limited"
$week.Range("D377").Value = "500480"
$week.Range("E377").Value = -1.61
$week.Range("F377").Value = 3623.5
$week.Range("G377").Value = 486960
$week.Range("H377").Value = "week"
$week.Range("I377").Value = "04/10/2024 11:32:49"

$week.Range("A378").Value = 5
$week.Range("B378").Value = "TORNTPHARM"
$week.Range("C378").Value = "Torrent Pharmaceuticals Limited"
$week.Range("D378").Value = "500420"
$week.Range("E378").Value = 2.4
$week.Range("F378").Value = 3473.55
$week.Range("G378").Value = 475367
$week.Range("H378").Value = "week"
$week.Range("I378").Value = "04/10/2024 11:32:49"

$week.Range("A379").Value = 6
$week.Range("B379").Value = "RELIANCE"
$week.Range("C379").Value = "Reliance Industries Limited"
$week.Range("D379").Value = "500325"
$week.Range("E379").Value = -1.45
$week.Range("F379").Value = 2773.05
$week.Range("G379").Value = 18536438
$week.Range("H379").Value = "week"
$week.Range("I379").Value = "04/10/2024 11:32:49"

$week.Range("A380").Value = 7
$week.Range("B380").Value = "GRANULES"
$week.Range("C380").Value = "Granules India Limited"
$week.Range("D380").Value = "532482"
$week.Range("E380").Value = -4.32
$week.Range("F380").Value = 569.75
$week.Range("G380").Value = 5724623
$week.Range("H380").Value = "week"
$week.Range("I380").Value = "04/10/2024 11:32:49"

$week.Range("A381").Value = 8
$week.Range("B381").Value = "COALINDIA"
$week.Range("C381").Value = "Coal India Limited"
$week.Range("D381").Value = "533278"
$week.Range("E381").Value = -1.03
$week.Range("F381").Value = 497.2
$week.Range("G381").Value = 7170285
$week.Range("H381").Value = "week"
$week.Range("I381").Value = "04/10/2024 11:32:49"

$week.Range("A382").Value = 9
$week.Range("B382").Value = "INDUSTOWER"
$week.Range("C382").Value = "Indus Towers Ltd (Bharti Infratel)"
$week.Range("D382").Value = "534816"
$week.Range("E382").Value = -1.43
$week.Range("F382").Value = 372.2
$week.Range("G382").Value = 10290941
$week.Range("H382").Value = "week"
$week.Range("I382").Value = "04/10/2024 11:32:49"

$week.Range("A383").Value = 10
$week.Range("B383").Value = "ASHOKLEY"
$week.Range("C383").Value = "Ashok Leyland Limited"
$week.Range("D383").Value = "500477"
$week.Range("E383").Value = -2.3
$week.Range("F383").Value = 225.39
$week.Range("G383").Value = 17371289
$week.Range("H383").Value = "week"
$week.Range("I383").Value = "04/10/2024 11:32:49"

$week.Range("A384").Value = 11
$week.Range("B384").Value = "MANAPPURAM"
$week.Range("C384").Value = "Manappuram Finance Limited"
$week.Range("D384").Value = "531213"
$week.Range("E384").Value = -1.01
$week.Range("F384").Value = 189.3
$week.Range("G384").Value = 8165417
$week.Range("H384").Value = "week"
$week.Range("I384").Value = "04/10/2024 11:32:49"

$week.Range("A385").Value = 12
$week.Range("B385").Value = "IDEA"
$week.Range("C385").Value = "Idea Cellular Limited"
$week.Range("D385").Value = "532822"
$week.Range("E385").Value = -0.8100000000000001
$week.Range("F385").Value = 9.789999999999999
$week.Range("G385").Value = 453171770
$week.Range("H385").Value = "week"
$week.Range("I385").Value = "04/10/2024 11:32:49"
